$d = $word.ActiveDocument

$d.Content.Find.Execute("Ativação: 01/01/1996", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2022", 2)

$d.Content.Find.Execute("1062721 - José Benedito Marcomini", $true, $false, $false, $false, $false,
                         $true, 1, $false, "519033 - Carlos Yujiro Shigue", 2)

$d.Content.Find.Execute("Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes.", 2)

$d.Content.Find.Execute("A média do semestre será computada com base na relação:M=(P1+2P2)/3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A média do semestre será computada com base na relação: M=(A1+A2)/2", 2)

$d.Content.Find.Execute("A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Não cabe recuperação.", 2)
